$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format cells whose new price text would otherwise be parsed as a number,
# so they stay text (matching the existing inline-string "Price" column).
$textCells = @("D5", "D6", "D7", "D9", "D10", "D12", "D13", "D15", "D17", "D20", "D21", "D23", "D24", "D25", "D26", "D27", "D29", "D30", "D32", "D33", "D34", "D35", "D36", "D38", "D40", "D41", "D43", "D44", "D48", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "51.797.65"
$ws.Range("E2").Value = "  -0.33%  "

$ws.Range("D3").Value = "2.972.52"
$ws.Range("E3").Value = "  +1.47%  "

$ws.Range("D5").Value = "353.48"
$ws.Range("E5").Value = "  -1.02%  "

$ws.Range("D6").Value = "106.40"
$ws.Range("E6").Value = "  -3.94%  "

$ws.Range("D7").Value = "0.552"
$ws.Range("E7").Value = "  -3.07%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").Value = "0.601"
$ws.Range("E9").Value = "  -4.42%  "

$ws.Range("D10").Value = "37.58"
$ws.Range("E10").Value = "  -4.87%  "

$ws.Range("E11").Value = "  +2.54%  "

$ws.Range("D12").Value = "0.0848"
$ws.Range("E12").Value = "  -4.18%  "

$ws.Range("D13").Value = "18.90"
$ws.Range("E13").Value = "  -4.41%  "

$ws.Range("D14").Value = "3.446.73"
$ws.Range("E14").Value = "  +1.69%  "

$ws.Range("D15").Value = "7.50"
$ws.Range("E15").Value = "  -5.10%  "

$ws.Range("D16").Value = "2.963.30"
$ws.Range("E16").Value = "  +0.59%  "

$ws.Range("D17").Value = "0.990"
$ws.Range("E17").Value = "  +0.15%  "

$ws.Range("D18").Value = "51.733.34"
$ws.Range("E18").Value = "  -0.48%  "

$ws.Range("E19").Value = "  +0.37%  "

$ws.Range("D20").Value = "7.36"
$ws.Range("E20").Value = "  -3.08%  "

$ws.Range("D21").Value = "13.33"
$ws.Range("E21").Value = "  -5.60%  "

$ws.Range("D22").Value = "0.0₃0960"
$ws.Range("E22").Value = "  -2.46%  "

$ws.Range("D23").Value = "68.80"
$ws.Range("E23").Value = "  -3.10%  "

$ws.Range("D24").Value = "262.34"
$ws.Range("E24").Value = "  -3.18%  "

$ws.Range("D25").Value = "2.69"
$ws.Range("E25").Value = "  -4.93%  "

$ws.Range("D26").Value = "0.175"
$ws.Range("E26").Value = "  -4.35%  "

$ws.Range("D27").Value = "26.66"
$ws.Range("E27").Value = "  -2.13%  "

$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("D29").Value = "7.31"
$ws.Range("E29").Value = "  -2.32%  "

$ws.Range("D30").Value = "0.110"
$ws.Range("E30").Value = "  +1.94%  "

$ws.Range("E31").Value = "  +2.86%  "

$ws.Range("D32").Value = "10.06"
$ws.Range("E32").Value = "  -5.37%  "

$ws.Range("B33").Value = "Toncoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D33").Value = "2.16"
$ws.Range("E33").Value = "  +12.27%  "

$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").Value = "35.54"
$ws.Range("E34").Value = "  -7.89%  "

$ws.Range("D35").Value = "50.89"
$ws.Range("E35").Value = "  -2.74%  "

$ws.Range("D36").Value = "0.0427"
$ws.Range("E36").Value = "  -4.01%  "

$ws.Range("E37").Value = "  +0.06%  "

$ws.Range("D38").Value = "3.22"
$ws.Range("E38").Value = "  -0.96%  "

$ws.Range("E39").Value = "  +0.61%  "

$ws.Range("D40").Value = "1.92"
$ws.Range("E40").Value = "  -4.47%  "

$ws.Range("D41").Value = "17.31"
$ws.Range("E41").Value = "  -6.65%  "

$ws.Range("E42").Value = "  -3.88%  "

$ws.Range("D43").Value = "22.83"
$ws.Range("E43").Value = "  -1.03%  "

$ws.Range("D44").Value = "123.54"
$ws.Range("E44").Value = "  +3.78%  "

$ws.Range("E45").Value = "  -0.32%  "

$ws.Range("D46").Value = "2.106.29"
$ws.Range("E46").Value = "  -1.67%  "

$ws.Range("E47").Value = "  -5.54%  "

$ws.Range("D48").Value = "2.32"
$ws.Range("E48").Value = "  -7.71%  "

$ws.Range("E49").Value = "  -4.19%  "

$ws.Range("D50").Value = "0.0328"
$ws.Range("E50").Value = "  -1.68%  "

$ws.Range("B51").Value = "SEI"
$ws.Range("C51").Value = "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
$ws.Range("D51").Value = "0.892"
$ws.Range("E51").Value = "  -2.38%  "
